$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas to round the result to 0 decimal places
$ws.Range("G4").Formula = "=ROUND(MOD(SQRT(G1^3+4*G1+1),11),0)"
$ws.Range("I4").Formula = "=ROUND(MOD(SQRT(I1^3+4*I1+1),11),0)"
$ws.Range("G5").Formula = "=ROUND(MOD(-SQRT(G1^3+4*G1+1),11),0)"
$ws.Range("I5").Formula = "=ROUND(MOD(-SQRT(I1^3+4*I1+1),11),0)"

# Align the style of A4/A5 with A1/A2/A3 (drop the extra applyFill variant)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the view: zoom and selection
$ws.Select()
$ws.Range("A1:L5").Select()
$excel.ActiveWindow.Zoom = 183
